# Add CO2_flux to the gap filling configuration (Berge_MDS sheet, cell C5)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Berge_MDS")
$ws.Range("C5").Value = "CO2_flux"
